$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Fill in row 34 with the new results entry
$ws.Range("A34").Value = "Advanced"
$ws.Range("B34").Value = "1-26,28-32,40-53"
$ws.Range("C34").Value = "20"
$ws.Range("D34").Value = "Yes"
$ws.Range("E34").Value = "onetree"
$ws.Range("F34").Value = "80.1"
$ws.Range("G34").Value = "1:49:33"
$ws.Range("H34").Value = "1:16"
$ws.Range("I34").Value = "1:03"
$ws.Range("J34").Value = "4 i7 CPUs, 16 GRAM"
$ws.Range("K34").Value = "results_27_01-154450"

# Match the "Machine" column style used by other data rows (centered, general format)
$ws.Range("J34").Style = "Normal"
$ws.Range("J34").HorizontalAlignment = -4108

# Update selection to the newly edited cell
$ws.Range("H34").Select()
